# COVAC tracker metadata reference workbook — update to v1.1.1
# (DHIS2 2.33 / 2.34 / 2.35 release)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "Package info" sheet — bump version / timestamps / identifier
# ---------------------------------------------------------------------
$wsInfo = $wb.Worksheets.Item("Package info")
$wsInfo.Range("B4").Value = "V1.1.1"
$wsInfo.Range("B6").Value = "20210310T015409"
$wsInfo.Range("B7").Value = "COVAC_TRACKER_V1.1.1_DHIS2.35.1-bca7d7b_20210310T015409"

# ---------------------------------------------------------------------
# 2. "dataElements" sheet — rename the code of "COVAC - AEFIs present"
# ---------------------------------------------------------------------
$wsDE = $wb.Worksheets.Item("dataElements")
$wsDE.Range("C2").Value = "COVAC_AEFIs_present"

# ---------------------------------------------------------------------
# 3. "dataElementGroups" sheet — swap the two data elements listed for
#    the "COVAC - Covid-19 vaccination registry" group (rows 16 & 17)
# ---------------------------------------------------------------------
$wsDEG = $wb.Worksheets.Item("dataElementGroups")
$wsDEG.Range("B16").Value = "COVAC - AEFIs present"
$wsDEG.Range("B17").Value = "COVAC - Cardiovascular Disease"

# ---------------------------------------------------------------------
# 4. "optionSets" sheet — fix typo in trimester option set description
# ---------------------------------------------------------------------
$wsOS = $wb.Worksheets.Item("optionSets")
$wsOS.Range("D5").Value = "1st Trimester (1-12 weeks); 2nd Trimester (13-28 weeks); 3rd Trimester (29-40 weeks)"

# ---------------------------------------------------------------------
# 5. "options" sheet — fix typo in "3rd Trimester" option name
# ---------------------------------------------------------------------
$wsOpt = $wb.Worksheets.Item("options")
$wsOpt.Range("B7").Value = "3rd Trimester (29-40 weeks)"

# ---------------------------------------------------------------------
# 6. "visualizations" sheet — rows 2-12 were reordered
# ---------------------------------------------------------------------
$wsViz = $wb.Worksheets.Item("visualizations")
$wsViz.Range("A2").Value = "COVAC - People with completed vaccination schedule"
$wsViz.Range("B2").Value = "COVAC - People with completed vaccination schedule"
$wsViz.Range("D2").Value = "TWG0cq8P539"

$wsViz.Range("A3").Value = "COVAC - People receiving COV-2 vs People completing the vaccination schedule (Cov-C)"
$wsViz.Range("B3").Value = ""
$wsViz.Range("D3").Value = "wHd33PaphEC"

$wsViz.Range("A4").Value = "COVAC - At least one underlying condition"
$wsViz.Range("B4").Value = ""
$wsViz.Range("D4").Value = "gNsB9zivLTy"

$wsViz.Range("A5").Value = "COVAC - Underlying conditions"
$wsViz.Range("B5").Value = ""
$wsViz.Range("D5").Value = "vFkbMQiABfj"

$wsViz.Range("A6").Value = "COVAC - Number of doses administered"
$wsViz.Range("B6").Value = ""
$wsViz.Range("D6").Value = "GJPPSQuVt4N"

$wsViz.Range("A7").Value = "COVAC - Complete vaccination uptake"
$wsViz.Range("B7").Value = ""
$wsViz.Range("D7").Value = "aUjo2Myd25f"

$wsViz.Range("A8").Value = "COVAC - Vaccine uptake by sex"
$wsViz.Range("B8").Value = ""
$wsViz.Range("D8").Value = "KV7fffdXnlY"

$wsViz.Range("A9").Value = "COVAC - Vaccine uptake by age group"
$wsViz.Range("B9").Value = ""
$wsViz.Range("D9").Value = "BWlYGFBDbO2"

$wsViz.Range("A10").Value = "COVAC - Vaccine uptake, last month"
$wsViz.Range("B10").Value = ""
$wsViz.Range("D10").Value = "MzSAvoJ0vLr"

$wsViz.Range("A11").Value = "COVAC - Drop-out from COV-1 to COV-c"
$wsViz.Range("B11").Value = ""
$wsViz.Range("D11").Value = "Hbs3xGj7XoN"

$wsViz.Range("A12").Value = "COVAC - Cumulative number of given doses, Last 4 weeks"
$wsViz.Range("B12").Value = "TEST TEST"
$wsViz.Range("D12").Value = "vmNUVdhuxN7"

# ---------------------------------------------------------------------
# 7. "trackedEntityAttributes" sheet — drop the "covac_" prefix on the
#    Sex attribute's code
# ---------------------------------------------------------------------
$wsTEA = $wb.Worksheets.Item("trackedEntityAttributes")
$wsTEA.Range("B11").Value = "patinfo_sex"

# ---------------------------------------------------------------------
# 8. "programs" sheet — bump last-updated date (kept as literal text,
#    not auto-converted to a date serial: force text format before the
#    write, then restore the original formatting from a neighbour cell)
# ---------------------------------------------------------------------
$wsProg = $wb.Worksheets.Item("programs")
$wsProg.Range("C2").NumberFormat = "@"
$wsProg.Range("C2").Value = "2021-03-08"
$wsProg.Range("A2").Copy()
$wsProg.Range("C2").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 9. "programRules" sheet — a new rule was added, shifting rows 21-55
#    down to 22-56 (dimension grows from A1:E55 to A1:E56)
# ---------------------------------------------------------------------
$wsPR = $wb.Worksheets.Item("programRules")
$wsPR.Rows.Item(21).Insert()

# new row inherits banding/format from the row 2 above it (even style);
# copy formatting from an "odd" row (e.g. row 19) so the zebra-stripe
# pattern continues correctly after the insert
$wsPR.Range("A19:E19").Copy()
$wsPR.Range("A21:E21").PasteSpecial(-4122)

$wsPR.Range("A21").Value = "R1bzqObecyQ"
$wsPR.Range("B21").Value = "Hide Suggested date for next dose if vaccine product has no more doses"
$wsPR.Range("C21").Value = "All vaccine types with two doses, after they receive one does, the ""next dose date"" will be hidden."
$wsPR.Range("D21").Value = ""
$wsPR.Range("E21").Value = "yDuAzyqYABS"
